$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade data rows (rows 10-12), matching columns A:H of existing data.
$data = @(
    @(8914.68,           8995.64,           18.84, 19.010000000000002, $true,  0.9,    42613.765543981484, $false),
    @(8963.7099999999991, 8914.68,           18.93, 18.824999999999999, $true, -0.55000000000000004, 42614.67287037037, $true),
    @(9011.2199999999993, 8963.7099999999991, 18.72, 18.62,             $true, -0.53,  42615.750150462962, $true)
)

$row = 10
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 7).NumberFormat = "m/d/yy h:mm"
    $ws.Cells.Item($row, 8).Value = $r[7]
    $row++
}
